# Updates Case_7_1 res_line loading_percent.xlsx sheet1:
# - extends header row with columns P (14) and Q (15), copying O1's style
# - recomputes B:H, removes column I entirely, resets O to 0, and adds
#   new P (0) / Q (computed) columns for every data row (2-25)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): add P1 = 14, Q1 = 15, matching O1's style ---
$ws.Range("O1").Copy() | Out-Null
$ws.Range("P1:Q1").PasteSpecial(-4122) | Out-Null
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# --- Data rows (2-25): columns B,C,D,E,F,G,H,J,K,L,M,N,O,P,Q ---
$colOrder = @("B","C","D","E","F","G","H","J","K","L","M","N","O","P","Q")

$rowsData = @"
24.89789934150223,18.25301277263919,4.054982730462291,29.63483041703122,18.13056958308222,2.060654926766087,3.416461098028841,0,0,0,0,0,0,0,13.93177726867801
23.26383590219017,17.1409390618791,3.927562455213644,27.71605081910991,17.49729497592208,2.066272940429423,3.139479544576587,0,0,0,0,0,0,0,13.6834973561419
22.20182512873893,16.42130942104709,3.846509554043123,26.47378002708209,17.11988626078928,2.069814707263204,2.964060937160555,0,0,0,0,0,0,0,13.54704538661909
21.75423156196883,16.11886727120157,3.812799646780699,25.95137399880545,16.96911306271626,2.071281952138181,2.891197615986665,0,0,0,0,0,0,0,13.49535625111767
21.67901269155977,16.0680946933409,3.807161944071675,25.86365133915934,16.94426426612312,2.071527054272306,2.879014494230377,0,0,0,0,0,0,0,13.48700641639094
22.19584878196305,16.41726766669485,3.846057645086946,26.46680019416365,17.1178404397331,2.069834397130231,2.96308390270716,0,0,0,0,0,0,0,13.54633258447055
24.34661721556397,17.87724318553504,4.011648518784196,28.98651297185654,17.90993640170255,2.06257332276784,3.322047863456711,0,0,0,0,0,0,0,13.84279295996338
28.10171957369256,20.44645016737692,4.312964283975089,33.42183000656939,19.54737350172488,2.049030941421981,3.985111876078617,0,0,0,0,0,0,0,14.55489348884665
30.57739167837056,22.14797249273774,4.557700116769004,35.49942372999039,20.59918098157599,2.039718306535877,4.410775765986231,0,0,0,0,0,0,0,14.99896247243963
31.59574558622074,22.83519410776042,5.007803317783687,28.94211995777329,19.4116084638581,2.037651252785357,4.74935085956407,0,0,0,0,0,0,0,13.81073806269926
31.95664881240601,23.07453827752758,5.347863969271172,23.04101367629282,18.16350039559847,2.037583979499019,5.512150096706766,0,0,0,0,0,0,0,12.69154346196812
31.84600825943449,22.99079515537812,5.627141582212349,17.19932329987481,16.75543168490039,2.039051616875743,6.489369618090509,0,0,0,0,0,0,0,11.53861850959021
31.55480574437501,22.78566660348587,5.797597134779451,13.22921809736451,15.68287818672621,2.040741287601649,7.26260981664733,0,0,0,0,0,0,0,10.71503463387819
31.38129251788887,22.66504992003988,5.825273799766417,12.28220740644145,15.38751297992019,2.041516563898398,7.440115290114314,0,0,0,0,0,0,0,10.51085250331039
30.38855116997392,21.98072700982351,5.676611318848624,12.1128561499128,15.13385129662749,2.045162024807705,7.152628156750532,0,0,0,0,0,0,0,10.51719187915193
29.7717255296767,21.55837638433653,5.46654749890296,14.11760469785684,15.53668309839994,2.047088132447558,6.468962159158048,0,0,0,0,0,0,0,10.96862797586405
29.42455603070891,21.32538497787987,5.184146049627055,18.57435640394323,16.53818868200793,2.047612175479546,5.459729313761628,0,0,0,0,0,0,0,11.846005298716
29.32904302826502,21.26956208387394,4.87199467733732,24.78989973588908,17.9288655093456,2.04677079658162,4.528639891267542,0,0,0,0,0,0,0,13.00025367589964
29.94551859568076,21.71300583584856,4.496327547852867,34.92596129223055,20.31058758239368,2.042152209980016,4.297821795127231,0,0,0,0,0,0,0,14.86886659007179
31.76864487013778,22.97154746301271,4.618726059219144,37.79561574836795,21.43184923100538,2.034696757944495,4.679217316707213,0,0,0,0,0,0,0,15.48986414217357
32.90832623376996,23.75852544307827,4.715552084575854,39.16346842533085,22.07066333072222,2.03000899010797,4.904705695923207,0,0,0,0,0,0,0,15.82816972328746
32.30447887015833,23.34145011795574,4.66415721179145,38.43822589831207,21.72916903209495,2.032506784222057,4.784645236438124,0,0,0,0,0,0,0,15.64612935729463
29.91974643493214,21.6968387252622,4.463354208521735,35.58510515996978,20.44822497922178,2.042064890147257,4.323337309462734,0,0,0,0,0,0,0,14.98899468366936
27.1349020761604,19.78301228928547,4.23409127005464,32.27560715676275,19.0962666902096,2.052625218022758,3.80987537854725,0,0,0,0,0,0,0,14.34776845306689
"@ -split "`n"

$r = 2
foreach ($line in $rowsData) {
    $vals = $line -split ","
    for ($i = 0; $i -lt $colOrder.Length; $i++) {
        $col = $colOrder[$i]
        $addr = "$col$r"
        $ws.Range($addr).Value = [double]$vals[$i]
    }
    # Column I is removed entirely for this row
    $ws.Range("I$r").ClearContents()
    $r++
}
